$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.184.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.568.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.19"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.569.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.53%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  +3.55%  "

$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.51%  "

$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.172.64"
$ws.Range("D13").Style = "Normal"

$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("E15").Value = "  -0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.592.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.42%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.280.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("E18").Value = "  -0.95%  "

$ws.Range("E19").Value = "  +8.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "

$ws.Range("E21").Value = "  +1.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "429.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.614"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.709.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.41%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +3.20%  "

$ws.Range("E28").Value = "  +2.39%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  -3.49%  "

$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("E32").Value = "  -0.89%  "

$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.560.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.46%  "

$ws.Range("E35").Value = "  -6.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.88"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "175.18"
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = "  -2.07%  "

$ws.Range("E43").Value = "  +2.13%  "

$ws.Range("E44").Value = "  +0.73%  "

$ws.Range("E45").Value = "  +1.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.11%  "

$ws.Range("E47").Value = "  +0.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.82%  "

$ws.Range("E49").Value = "  +2.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.09%  "
